$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.583.50'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.752.39'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4483'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.58%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3564'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07467'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.081'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.71'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.983'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.144'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.756.38'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.58'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001055'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06391'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.730'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.637.72'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("E24").Value = '  -0.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.090'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.69'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.958.02'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.088'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.55'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.090'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09175'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.664'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.505'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02285'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.06%  '
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2094'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.62%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06022'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6285'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.924'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.63%  '
$ws.Range("E41").Value = '  -0.48%  '
$ws.Range("E42").Value = '  -0.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.755'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.12'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.04%  '
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5862'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.936'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06886'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.128'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.92%  '
